$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.648.02'
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').Value = '1.603.68'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.517'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '27.14'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '43.42'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.53%  '
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D13').Value = '1.833.88'
$ws.Range('E13').Value = '  +2.60%  '
$ws.Range('D14').Value = '1.600.78'
$ws.Range('E14').Value = '  +2.33%  '
$ws.Range('D15').Value = '29.660.18'
$ws.Range('E15').Value = '  +3.46%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.536'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.83%  '
$ws.Range('E17').Value = '  +2.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.44'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('E21').Value = '  +1.20%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.22'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.93%  '
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.79%  '
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0474'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.58%  '
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.13'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.33%  '
$ws.Range('B35').Value = 'Maker'
$ws.Range('C35').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D35').Value = '1.431.56'
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('E37').Value = '  +4.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.82'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('E40').Value = '  +1.38%  '
$ws.Range('E41').Value = '  +3.81%  '
$ws.Range('E42').Value = '  +1.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '54.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +29.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0488'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.800'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.25%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '65.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.954'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +13.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.28'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').Value = '1.744.24'
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '86.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.03%  '
